# --- jobs sheet: Job Release Date column removed; Job Due Date now a real datetime ---
$wb = $excel.ActiveWorkbook
$jobs = $wb.Worksheets.Item("jobs")

# Remove the "Job Release Date" column (old column C, all zeros)
$jobs.Columns.Item(3).Delete()

# Column C now holds "Job Due Date" (previously minutes-since-epoch in column D).
# Replace with real datetimes and format as date+time.
$jobs.Range("C2:C11").NumberFormat = "m/d/yy h:mm"
$jobs.Range("C2").Value = 43739.75
$jobs.Range("C3").Value = 43740.75
$jobs.Range("C4").Value = 43740.75
$jobs.Range("C5").Value = 43740.75
$jobs.Range("C6").Value = 43739.75
$jobs.Range("C7").Value = 43742.75
$jobs.Range("C8").Value = 43741.75
$jobs.Range("C9").Value = 43740.75
$jobs.Range("C10").Value = 43739.75
$jobs.Range("C11").Value = 43739.75

$jobs.PageSetup.Orientation = 1
$jobs.Range("C1").Select()

# --- machines sheet: widen the first availability-window column ---
$machines = $wb.Worksheets.Item("machines")
$machines.Columns.Item(3).ColumnWidth = 46.109375

# --- tasks sheet: reset scroll position ---
$tasks = $wb.Worksheets.Item("tasks")
$tasks.Range("C47:C51").Select()

# Restore original active sheet/selection (machines tab was active before editing)
$machines.Activate()
$machines.Range("H8").Select()
